$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10637.714
$ws.Range("I32").Value = 10562.333
$ws.Range("J32").Value = 10694.25
$ws.Range("K32").Value = 10562.333
$ws.Range("L32").Value = 10694.25
$ws.Range("M32").Value = -10236.333
$ws.Range("N32").Value = -11346.25
$ws.Range("H62").Value = 7119.6
$ws.Range("I62").Value = 6271.0586
$ws.Range("J62").Value = 8229.23
$ws.Range("K62").Value = 6271.0586
$ws.Range("L62").Value = 8229.23
$ws.Range("M62").Value = -5647.0586
$ws.Range("N62").Value = -9477.23
$ws.Range("H65").Value = 7119.6
$ws.Range("I65").Value = 6271.0586
$ws.Range("J65").Value = 8229.23
$ws.Range("K65").Value = 31355.293
$ws.Range("L65").Value = 41146.14999999999
$ws.Range("M65").Value = -28235.293
$ws.Range("N65").Value = -47386.14999999999
$ws.Range("H70").Value = 14361396
$ws.Range("I70").Value = 202899.6
$ws.Range("K70").Value = 608698.8
$ws.Range("M70").Value = -608428.8
$ws.Range("H73").Value = 14361396
$ws.Range("I73").Value = 202899.6
$ws.Range("K73").Value = 608698.8
$ws.Range("M73").Value = -607762.8
$ws.Range("H100").Value = 5036.3335
$ws.Range("I100").Value = 2968.4
$ws.Range("J100").Value = 7621.25
$ws.Range("K100").Value = 2968.4
$ws.Range("L100").Value = 7621.25
$ws.Range("M100").Value = -2427.4
$ws.Range("N100").Value = -8703.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3766.8333
$ws.Range("I2").Value = 702.5484
$ws.Range("K2").Value = 702.5484
$ws.Range("M2").Value = -589.5484
$ws.Range("H45").Value = 2710.1538
$ws.Range("I45").Value = 2212.3809
$ws.Range("J45").Value = 4800.8
$ws.Range("K45").Value = 2212.3809
$ws.Range("L45").Value = 4800.8
$ws.Range("M45").Value = -1835.3809
$ws.Range("N45").Value = -5554.8
$ws.Range("H74").Value = 12347534
$ws.Range("I74").Value = 15152530
$ws.Range("J74").Value = 5551.8
$ws.Range("K74").Value = 15152530
$ws.Range("L74").Value = 5551.8
$ws.Range("M74").Value = -15151656
$ws.Range("N74").Value = -7299.8
$ws.Range("H77").Value = 12347534
$ws.Range("I77").Value = 15152530
$ws.Range("J77").Value = 5551.8
$ws.Range("K77").Value = 75762650
$ws.Range("L77").Value = 27759
$ws.Range("M77").Value = -75758282
$ws.Range("N77").Value = -36495
$ws.Range("H116").Value = 3766.8333
$ws.Range("I116").Value = 702.5484
$ws.Range("K116").Value = 702.5484
$ws.Range("M116").Value = 1591.4516
$ws.Range("H132").Value = 2806.2163
$ws.Range("I132").Value = 2159.9285
$ws.Range("J132").Value = 4816.8887
$ws.Range("K132").Value = 6479.7855
$ws.Range("L132").Value = 14450.6661
$ws.Range("M132").Value = -3949.7855
$ws.Range("N132").Value = -19510.6661

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3766.8333
$ws.Range("I3").Value = 702.5484
$ws.Range("K3").Value = 702.5484
$ws.Range("M3").Value = -588.5484
$ws.Range("H86").Value = 5649.826
$ws.Range("I86").Value = 4821.95
$ws.Range("K86").Value = 4821.95
$ws.Range("M86").Value = -3698.95
$ws.Range("H89").Value = 5649.826
$ws.Range("I89").Value = 4821.95
$ws.Range("K89").Value = 24109.75
$ws.Range("M89").Value = -18493.75
$ws.Range("H105").Value = 15850.917
$ws.Range("I105").Value = 15446.467
$ws.Range("K105").Value = 15446.467
$ws.Range("M105").Value = -13699.467
$ws.Range("H134").Value = 2511.8
$ws.Range("I134").Value = 1671.9
$ws.Range("K134").Value = 5015.700000000001
$ws.Range("M134").Value = -2480.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26355.795
$ws.Range("I31").Value = 2350.879
$ws.Range("J31").Value = 98370.55
$ws.Range("K31").Value = 2350.879
$ws.Range("L31").Value = 98370.55
$ws.Range("M31").Value = -2055.879
$ws.Range("N31").Value = -98960.55
$ws.Range("H34").Value = 26355.795
$ws.Range("I34").Value = 2350.879
$ws.Range("J34").Value = 98370.55
$ws.Range("K34").Value = 2350.879
$ws.Range("L34").Value = 98370.55
$ws.Range("M34").Value = -2148.879
$ws.Range("N34").Value = -98774.55
$ws.Range("H58").Value = 3656.36
$ws.Range("I58").Value = 1483.8
$ws.Range("J58").Value = 6915.2
$ws.Range("K58").Value = 1483.8
$ws.Range("L58").Value = 6915.2
$ws.Range("M58").Value = -1280.8
$ws.Range("N58").Value = -7321.2
$ws.Range("H99").Value = 2281.125
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2749.6667
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2749.6667
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -5745.6667
$ws.Range("H126").Value = 2281.125
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2749.6667
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 8249.000100000001
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -13189.0001
$ws.Range("H134").Value = 2421.6365
$ws.Range("I134").Value = 1221.8125
$ws.Range("J134").Value = 5621.1665
$ws.Range("K134").Value = 3665.4375
$ws.Range("L134").Value = 16863.4995
$ws.Range("M134").Value = -1130.4375
$ws.Range("N134").Value = -21933.4995
$ws.Range("H136").Value = 3656.36
$ws.Range("I136").Value = 1483.8
$ws.Range("J136").Value = 6915.2
$ws.Range("K136").Value = 4451.4
$ws.Range("L136").Value = 20745.6
$ws.Range("M136").Value = -1901.4
$ws.Range("N136").Value = -25845.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 43030868
$ws.Range("I4").Value = 64117840
$ws.Range("K4").Value = 192353520
$ws.Range("M4").Value = -192353408
$ws.Range("H131").Value = 5809848
$ws.Range("J131").Value = 3670587.8
$ws.Range("L131").Value = 11011763.4
$ws.Range("N131").Value = -11021843.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7515.4
$ws.Range("I70").Value = 6327.6665
$ws.Range("J70").Value = 9297
$ws.Range("K70").Value = 6327.6665
$ws.Range("L70").Value = 9297
$ws.Range("M70").Value = -6057.6665
$ws.Range("N70").Value = -9837
$ws.Range("H73").Value = 7515.4
$ws.Range("I73").Value = 6327.6665
$ws.Range("J73").Value = 9297
$ws.Range("K73").Value = 6327.6665
$ws.Range("L73").Value = 9297
$ws.Range("M73").Value = -5391.6665
$ws.Range("N73").Value = -11169
$ws.Range("H113").Value = 2171.647
$ws.Range("I113").Value = 1851.6428
$ws.Range("K113").Value = 1851.6428
$ws.Range("M113").Value = 318.3571999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 75902.60000000001
$ws.Range("J36").Value = 75902.60000000001
$ws.Range("L36").Value = 75902.60000000001
$ws.Range("N36").Value = -77026.60000000001
$ws.Range("H55").Value = 1854117.2
$ws.Range("I55").Value = 2778877
$ws.Range("K55").Value = 2778877
$ws.Range("M55").Value = -2778704
$ws.Range("H61").Value = 2828.8696
$ws.Range("I61").Value = 1872.8948
$ws.Range("J61").Value = 7369.75
$ws.Range("K61").Value = 1872.8948
$ws.Range("L61").Value = 7369.75
$ws.Range("M61").Value = -1670.8948
$ws.Range("N61").Value = -7773.75
$ws.Range("H100").Value = 12937
$ws.Range("I100").Value = 4683.5
$ws.Range("J100").Value = 15000.375
$ws.Range("K100").Value = 4683.5
$ws.Range("L100").Value = 15000.375
$ws.Range("M100").Value = -4142.5
$ws.Range("N100").Value = -16082.375
$ws.Range("H113").Value = 2828.8696
$ws.Range("I113").Value = 1872.8948
$ws.Range("J113").Value = 7369.75
$ws.Range("K113").Value = 1872.8948
$ws.Range("L113").Value = 7369.75
$ws.Range("M113").Value = 297.1052
$ws.Range("N113").Value = -11709.75
$ws.Range("H136").Value = 3177.8044
$ws.Range("I136").Value = 2387.647
$ws.Range("J136").Value = 5416.5835
$ws.Range("K136").Value = 7162.941
$ws.Range("L136").Value = 16249.7505
$ws.Range("M136").Value = -4612.941
$ws.Range("N136").Value = -21349.7505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 923.7838
$ws.Range("I100").Value = 742.8276
$ws.Range("K100").Value = 1485.6552
$ws.Range("M100").Value = -944.6551999999999
$ws.Range("H113").Value = 323.62964
$ws.Range("I113").Value = 310.14285
$ws.Range("J113").Value = 370.83334
$ws.Range("K113").Value = 930.4285500000001
$ws.Range("L113").Value = 1112.50002
$ws.Range("M113").Value = 1239.57145
$ws.Range("N113").Value = -5452.500019999999
$ws.Range("H136").Value = 4990.16
$ws.Range("I136").Value = 3464.7144
$ws.Range("K136").Value = 10394.1432
$ws.Range("M136").Value = -7844.143199999999
